$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Values such as "1.210" or "0.3731" look numeric; force the target
# cells to text format first so Excel keeps the original literal text
# representation (matching the source inline-string cells) instead of
# silently converting them to a float.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "22.429.21"
$ws.Range("E2").Value = "  +0.10%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.571.30"
$ws.Range("E3").Value = "  -0.07%  "
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "1.004"
$ws.Range("E5").Value = "  +0.14%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "291.38"
$ws.Range("E6").Value = "  +0.21%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3731"
$ws.Range("E7").Value = "  -0.97%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "49.80"
$ws.Range("E8").Value = "  -0.13%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3387"
$ws.Range("E9").Value = "  -0.97%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07541"
$ws.Range("E10").Value = "  -1.48%  "
$ws.Range("E11").Value = "  -2.42%  "
$ws.Range("E12").Value = "  +0.04%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "21.34"
$ws.Range("E13").Value = "  +0.32%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.977"
$ws.Range("E14").Value = "  -0.60%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.921"
$ws.Range("E15").Value = "  -0.28%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.573.33"
$ws.Range("E16").Value = "  +0.16%  "
$ws.Range("E17").Value = "  -1.40%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "90.79"
$ws.Range("E18").Value = "  +0.90%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06735"
$ws.Range("E19").Value = "  -0.60%  "
$ws.Range("E20").Value = "  +0.13%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.256"
$ws.Range("E21").Value = "  +0.60%  "
$ws.Range("E22").Value = "  -3.08%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "12.11"
$ws.Range("E23").Value = "  +0.61%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "22.422.18"
$ws.Range("E24").Value = "  +0.05%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.334"
$ws.Range("E25").Value = "  -3.51%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.615"
$ws.Range("E26").Value = "  -3.61%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "20.08"
$ws.Range("E27").Value = "  -0.88%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "148.14"
$ws.Range("E28").Value = "  +1.07%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.015"
$ws.Range("E29").Value = "  -0.46%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "125.38"
$ws.Range("E30").Value = "  -0.72%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.748.19"
$ws.Range("E31").Value = "  +0.14%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.042"
$ws.Range("E32").Value = "  +5.72%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.111"
$ws.Range("E33").Value = "  -1.30%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.969"
$ws.Range("E34").Value = "  -1.98%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "9.719"
$ws.Range("E35").Value = "  -3.18%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.08349"
$ws.Range("E36").Value = "  -2.94%  "
$ws.Range("E37").Value = "  +3.76%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02462"
$ws.Range("E38").Value = "  -3.36%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.2283"
$ws.Range("E39").Value = "  -1.41%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.06501"
$ws.Range("E40").Value = "  -1.16%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.433"
$ws.Range("E41").Value = "  -0.73%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "11.22"
$ws.Range("E42").Value = "  -2.90%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.6192"
$ws.Range("E43").Value = "  -3.89%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.003"
$ws.Range("E44").Value = "  +0.10%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.85"
$ws.Range("E45").Value = "  -2.01%  "
$ws.Range("E46").Value = "  +0.50%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5779"
$ws.Range("E47").Value = "  -3.86%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "129.47"
$ws.Range("E48").Value = "  +3.13%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.067"
$ws.Range("E49").Value = "  -0.82%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.210"
$ws.Range("E50").Value = "  -6.99%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.07308"
$ws.Range("E51").Value = "  -0.24%  "
